# Update the "想去人数" (F column) counts for specific events in both the
# "展览" sheet and the "全部类型" sheet, per the commit's regenerated data.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 243
$wsExhibit.Range("F7").Value = 550
$wsExhibit.Range("F9").Value = 6770
$wsExhibit.Range("F13").Value = 175
$wsExhibit.Range("F15").Value = 1090
$wsExhibit.Range("F16").Value = 16151
$wsExhibit.Range("F17").Value = 1583
$wsExhibit.Range("F22").Value = 11325
$wsExhibit.Range("F24").Value = 956
$wsExhibit.Range("F25").Value = 4453
$wsExhibit.Range("F26").Value = 308

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 243
$wsAll.Range("F7").Value = 550
$wsAll.Range("F10").Value = 6770
$wsAll.Range("F14").Value = 175
$wsAll.Range("F17").Value = 1090
$wsAll.Range("F18").Value = 16151
$wsAll.Range("F19").Value = 1583
$wsAll.Range("F26").Value = 11325
$wsAll.Range("F28").Value = 956
$wsAll.Range("F29").Value = 4453
$wsAll.Range("F30").Value = 308
